$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'302.56"
$ws.Range("E2").Value = "'0.69%"
$ws.Range("D3").Value = "'31.89"
$ws.Range("E3").Value = "'0.37%"
$ws.Range("D4").Value = "'5.008"
$ws.Range("D5").Value = "'0.07880"
$ws.Range("D6").Value = "'2.081"
$ws.Range("E6").Value = "'-19.70%"
$ws.Range("D7").Value = "'7.823"
$ws.Range("E7").Value = "'-0.06%"
$ws.Range("D8").Value = "'3.795"
$ws.Range("E8").Value = "'-1.22%"
$ws.Range("D9").Value = "'0.9285"
$ws.Range("E9").Value = "'0.38%"
$ws.Range("D10").Value = "'0.1757"
$ws.Range("E10").Value = "'0.22%"
$ws.Range("D11").Value = "'0.07967"
$ws.Range("E11").Value = "'6.52%"
$ws.Range("D12").Value = "'0.08835"
$ws.Range("E12").Value = "'-0.75%"
$ws.Range("E13").Value = "'4.35%"
$ws.Range("E14").Value = "'0.32%"
$ws.Range("D15").Value = "'0.001511"
$ws.Range("E15").Value = "'0.36%"
$ws.Range("D16").Value = "'0.005977"
$ws.Range("E16").Value = "'3.16%"
$ws.Range("E17").Value = "'-3.49%"
$ws.Range("D18").Value = "'2.280"
$ws.Range("E18").Value = "'0.83%"
$ws.Range("D19").Value = "'0.3293"
$ws.Range("E20").Value = "'-3.98%"
$ws.Range("D21").Value = "'4.158"
$ws.Range("E21").Value = "'6.43%"
$ws.Range("D22").Value = "'0.1792"
$ws.Range("E22").Value = "'5.63%"
$ws.Range("D23").Value = "'0.04606"
$ws.Range("E23").Value = "'0.01%"
$ws.Range("D24").Value = "'0.001238"
$ws.Range("E24").Value = "'-0.40%"
$ws.Range("D25").Value = "'0.004501"
$ws.Range("E25").Value = "'-0.56%"
$ws.Range("D26").Value = "'0.0001250"
$ws.Range("E26").Value = "'4.36%"
$ws.Range("D39").Value = "'0.01735"
$ws.Range("E39").Value = "'-1.91%"
$ws.Range("D40").Value = "'0.04760"
$ws.Range("E40").Value = "'4.92%"
$ws.Range("D41").Value = "'0.007339"
$ws.Range("E41").Value = "'5.16%"
$ws.Range("D42").Value = "'0.1366"
$ws.Range("E42").Value = "'-0.83%"
$ws.Range("D43").Value = "'0.002340"
$ws.Range("E43").Value = "'9.47%"
$ws.Range("D44").Value = "'0.01095"
$ws.Range("E44").Value = "'13.95%"
$ws.Range("D45").Value = "'0.00006059"
$ws.Range("E45").Value = "'-5.11%"
$ws.Range("E46").Value = "'0.20%"
$ws.Range("D47").Value = "'0.003400"
$ws.Range("E47").Value = "'-61.10%"
$ws.Range("D48").Value = "'0.8205"
$ws.Range("E48").Value = "'2.06%"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("E49").Value = "'0.20%"
$ws.Range("D50").Value = "'0.0002000"
$ws.Range("E50").Value = "'0.20%"
